$wb = $excel.ActiveWorkbook

# --- Update panel_query_time values on the "data" sheet (column F) ---
$ws1 = $wb.Worksheets.Item("data")
$ws1.Range("F2").Value = "2021-10-05 14:33:57.161778"
$ws1.Range("F3").Value = "2021-10-05 14:33:57.161787"
$ws1.Range("F4").Value = "2021-10-05 14:33:57.161790"
$ws1.Range("F5").Value = "2021-10-05 14:33:57.161793"
$ws1.Range("F6").Value = "2021-10-05 14:33:57.161796"
$ws1.Range("F7").Value = "2021-10-05 14:33:57.161798"

# --- Add a new "metadata" worksheet right after the "data" sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Hereditary Haemorrhagic Telangiectasia"
$ws2.Range("C2").Value = 260
$ws2.Range("D2").Value = "'1.0"
$ws2.Range("D2").ClearFormats()
$ws2.Range("E2").Value = "2021-01-21T10:54:52.442904Z"
$ws2.Range("F2").Value = "2021-10-05 14:33:57.158467"
$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/260/?format=json"

# Apply the same header formatting (bold, bordered, centered) used on the
# "data" sheet's header row to the new sheet's header row and the A2 index cell.
$ws1.Range("B1").Copy()
$ws2.Range("B1:G1").PasteSpecial(-4122)
$ws2.Range("A2").PasteSpecial(-4122)
